# Auto-generated PowerShell COM-interop script replicating the commit diff.
# "Atualizando o arquivo XLSX" -- refresh of the weekly FlashScore odds sheet:
#   - 2 new fixtures inserted as rows 5 and 6 (EGYPT - PREMIER LEAGUE, POLAND - DIVISION 1)
#   - the 3 pre-existing Saudi fixtures shift down to rows 7-9
#   - their odds are refreshed to the latest values from the diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 5, pushing the existing rows 5-7 down to 7-9
$ws.Rows("5:6").Insert()

# Fill brand-new row 5: APE9ifU7 - ZED vs Al Ahly (EGYPT - PREMIER LEAGUE)
$ws.Range('A5').Value = 'APE9ifU7'
$ws.Range('B5').Value = '07/11/2024'
$ws.Range('C5').Value = '15:00'
$ws.Range('D5').Value = 'EGYPT - PREMIER LEAGUE'
$ws.Range('E5').Value = 'ZED'
$ws.Range('F5').Value = 'Al Ahly'
$ws.Range('G5').Value = 5.5
$ws.Range('H5').Value = 3.8
$ws.Range('I5').Value = 1.57
$ws.Range('J5').Value = 5.5
$ws.Range('K5').Value = 2.18
$ws.Range('L5').Value = 2.12
$ws.Range('M5').Value = 1.06
$ws.Range('N5').Value = 7.5
$ws.Range('O5').Value = 1.27
$ws.Range('P5').Value = 3.4
$ws.Range('Q5').Value = 1.82
$ws.Range('R5').Value = 1.93
$ws.Range('S5').Value = 1.39
$ws.Range('T5').Value = 2.77
$ws.Range('U5').Value = 1.83
$ws.Range('V5').Value = 1.87
$ws.Range('W5').Value = 14.5
$ws.Range('X5').Value = 32
$ws.Range('Y5').Value = 17
$ws.Range('Z5').Value = 110
$ws.Range('AA5').Value = 55
$ws.Range('AB5').Value = 55
$ws.Range('AC5').Value = 7.5
$ws.Range('AD5').Value = 7.3
$ws.Range('AE5').Value = 16.5
$ws.Range('AF5').Value = 75
$ws.Range('AG5').Value = 600
$ws.Range('AH5').Value = 6.8
$ws.Range('AI5').Value = 7.4
$ws.Range('AJ5').Value = 8
$ws.Range('AK5').Value = 11.5
$ws.Range('AL5').Value = 12.5
$ws.Range('AM5').Value = 25
$ws.Range('AN5').Value = 7
$ws.Range('AO5').Value = 32
$ws.Range('AP5').Value = 37
$ws.Range('AQ5').Value = 200
$ws.Range('AR5').Value = 250
$ws.Range('AS5').Value = 500
$ws.Range('AT5').Value = 2.77
$ws.Range('AU5').Value = 7.7
$ws.Range('AV5').Value = 75
$ws.Range('AW5').Value = 3.4
$ws.Range('AX5').Value = 7.6
$ws.Range('AY5').Value = 17.5
$ws.Range('AZ5').Value = 25
$ws.Range('BA5').Value = 55
$ws.Range('BB5').Value = 250
$ws.Range('BC5').Value = 51
$ws.Range('BD5').Value = 51

# Fill brand-new row 6: v12EbqCc - Ruch Chorzow vs Chrobry Glogow (POLAND - DIVISION 1)
$ws.Range('A6').Value = 'v12EbqCc'
$ws.Range('B6').Value = '07/11/2024'
$ws.Range('C6').Value = '15:00'
$ws.Range('D6').Value = 'POLAND - DIVISION 1'
$ws.Range('E6').Value = 'Ruch Chorzow'
$ws.Range('F6').Value = 'Chrobry Glogow'
$ws.Range('G6').Value = 1.57
$ws.Range('H6').Value = 3.8
$ws.Range('I6').Value = 5.5
$ws.Range('J6').Value = 2.1
$ws.Range('K6').Value = 2.3
$ws.Range('L6').Value = 5.5
$ws.Range('M6').Value = 1.04
$ws.Range('N6').Value = 13
$ws.Range('O6').Value = 1.25
$ws.Range('P6').Value = 3.75
$ws.Range('Q6').Value = 1.83
$ws.Range('R6').Value = 2.03
$ws.Range('S6').Value = 1.36
$ws.Range('T6').Value = 3
$ws.Range('U6').Value = 1.83
$ws.Range('V6').Value = 1.83
$ws.Range('W6').Value = 7.5
$ws.Range('X6').Value = 7.5
$ws.Range('Y6').Value = 8.5
$ws.Range('Z6').Value = 12
$ws.Range('AA6').Value = 13
$ws.Range('AB6').Value = 26
$ws.Range('AC6').Value = 11
$ws.Range('AD6').Value = 7.5
$ws.Range('AE6').Value = 17
$ws.Range('AF6').Value = 51
$ws.Range('AG6').Value = 251
$ws.Range('AH6').Value = 15
$ws.Range('AI6').Value = 29
$ws.Range('AJ6').Value = 17
$ws.Range('AK6').Value = 51
$ws.Range('AL6').Value = 41
$ws.Range('AM6').Value = 41
$ws.Range('AN6').Value = 3.6
$ws.Range('AO6').Value = 8
$ws.Range('AP6').Value = 19
$ws.Range('AQ6').Value = 23
$ws.Range('AR6').Value = 41
$ws.Range('AS6').Value = 126
$ws.Range('AT6').Value = 3
$ws.Range('AU6').Value = 8.5
$ws.Range('AV6').Value = 51
$ws.Range('AW6').Value = 7
$ws.Range('AX6').Value = 29
$ws.Range('AY6').Value = 34
$ws.Range('AZ6').Value = 101
$ws.Range('BA6').Value = 126
$ws.Range('BB6').Value = 251
$ws.Range('BC6').Value = 81
$ws.Range('BD6').Value = 81

# Refresh the odds that changed for the fixture now on row 7 (YyDJubM9 - Al Qadisiya vs Al Feiha)
$ws.Range('G7').Value = 1.3
$ws.Range('H7').Value = 5.25
$ws.Range('I7').Value = 7.5
$ws.Range('J7').Value = 1.73
$ws.Range('K7').Value = 2.5
$ws.Range('L7').Value = 7.5
$ws.Range('Q7').Value = 1.7
$ws.Range('R7').Value = 2.1
$ws.Range('U7').Value = 2.2
$ws.Range('V7').Value = 1.62
$ws.Range('X7').Value = 6
$ws.Range('Y7').Value = 9.5
$ws.Range('Z7').Value = 8
$ws.Range('AD7').Value = 11
$ws.Range('AF7').Value = 81
$ws.Range('AG7').Value = 1250
$ws.Range('AH7').Value = 19
$ws.Range('AJ7').Value = 23
$ws.Range('AK7').Value = 101
$ws.Range('AO7').Value = 6
$ws.Range('AQ7').Value = 15
$ws.Range('AU7').Value = 10
$ws.Range('AW7').Value = 9.5

# Refresh the odds that changed for the fixture now on row 8 (8fR1hy6F - Al Kholood vs Al Shabab)
$ws.Range('Q8').Value = 1.98
$ws.Range('R8').Value = 1.83
$ws.Range('AS8').Value = 400

# Row 9 (O6ibYFEq - Al Orubah vs Al Ittihad) is unchanged aside from the shift handled by Insert() above.

Write-Output 'edit complete'
